$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Cells.Item(32, 8).Value = 22223516  # H32 was 47621108
$ws_ALC.Cells.Item(32, 9).Value = 111111544  # I32 was 166667620
$ws_ALC.Cells.Item(32, 10).Value = 1508.5  # J32 was 2500.4
$ws_ALC.Cells.Item(32, 11).Value = 111111544  # K32 was 166667620
$ws_ALC.Cells.Item(32, 12).Value = 1508.5  # L32 was 2500.4
$ws_ALC.Cells.Item(32, 13).Value = -111111218  # M32 was -166667294
$ws_ALC.Cells.Item(32, 14).Value = -2160.5  # N32 was -3152.4
$ws_ALC.Cells.Item(33, 8).Value = 312.23077  # H33 was 265.7213
$ws_ALC.Cells.Item(33, 9).Value = 176.31818  # I33 was 210.17545
$ws_ALC.Cells.Item(33, 10).Value = 1059.75  # J33 was 1057.25
$ws_ALC.Cells.Item(33, 11).Value = 176.31818  # K33 was 210.17545
$ws_ALC.Cells.Item(33, 12).Value = 1059.75  # L33 was 1057.25
$ws_ALC.Cells.Item(33, 13).Value = 52.68181999999999  # M33 was 18.82454999999999
$ws_ALC.Cells.Item(33, 14).Value = -1517.75  # N33 was -1515.25
$ws_ALC.Cells.Item(41, 8).Value = 713.2  # H41 was 1216.6666
$ws_ALC.Cells.Item(41, 9).Value = 382.66666  # I41 was 950
$ws_ALC.Cells.Item(41, 10).Value = 933.55554  # J41 was 1350
$ws_ALC.Cells.Item(41, 11).Value = 382.66666  # K41 was 950
$ws_ALC.Cells.Item(41, 12).Value = 933.55554  # L41 was 1350
$ws_ALC.Cells.Item(41, 13).Value = 57.33334000000002  # M41 was -510
$ws_ALC.Cells.Item(41, 14).Value = -1813.55554  # N41 was -2230
$ws_ALC.Cells.Item(135, 8).Value = 68182630  # H135 was 93752640
$ws_ALC.Cells.Item(135, 9).Value = 47619900  # I135 was 90910530
$ws_ALC.Cells.Item(135, 10).Value = 500000000  # J135 was 100005300
$ws_ALC.Cells.Item(135, 11).Value = 428579100  # K135 was 818194770
$ws_ALC.Cells.Item(135, 12).Value = 4500000000  # L135 was 900047700
$ws_ALC.Cells.Item(135, 13).Value = -428576565  # M135 was -818192235
$ws_ALC.Cells.Item(135, 14).Value = -4500005070  # N135 was -900052770
$ws_ALC.Cells.Item(141, 8).Value = 5828  # H141 was 3982.7273
$ws_ALC.Cells.Item(141, 9).Value = 3775  # I141 was 2075
$ws_ALC.Cells.Item(141, 10).Value = 7001.143  # J141 was 9070
$ws_ALC.Cells.Item(141, 11).Value = 11325  # K141 was 6225
$ws_ALC.Cells.Item(141, 12).Value = 21003.429  # L141 was 27210
$ws_ALC.Cells.Item(141, 13).Value = -6145  # M141 was -1045
$ws_ALC.Cells.Item(141, 14).Value = -31363.429  # N141 was -37570
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Cells.Item(32, 8).Value = 12185.941  # H32 was 12941.4375
$ws_ARM.Cells.Item(32, 9).Value = 11335.512  # I32 was 12178.325
$ws_ARM.Cells.Item(32, 11).Value = 11335.512  # K32 was 12178.325
$ws_ARM.Cells.Item(32, 13).Value = -11048.512  # M32 was -11891.325
$ws_ARM.Cells.Item(63, 8).Value = 2722.5  # H63 was 2880.9375
$ws_ARM.Cells.Item(63, 9).Value = 1894.2307  # I63 was 2059.5833
$ws_ARM.Cells.Item(63, 10).Value = 4876  # J63 was 5345
$ws_ARM.Cells.Item(63, 11).Value = 1894.2307  # K63 was 2059.5833
$ws_ARM.Cells.Item(63, 12).Value = 4876  # L63 was 5345
$ws_ARM.Cells.Item(63, 13).Value = -1208.2307  # M63 was -1373.5833
$ws_ARM.Cells.Item(63, 14).Value = -6248  # N63 was -6717
$ws_ARM.Cells.Item(66, 8).Value = 2722.5  # H66 was 2880.9375
$ws_ARM.Cells.Item(66, 9).Value = 1894.2307  # I66 was 2059.5833
$ws_ARM.Cells.Item(66, 10).Value = 4876  # J66 was 5345
$ws_ARM.Cells.Item(66, 11).Value = 9471.1535  # K66 was 10297.9165
$ws_ARM.Cells.Item(66, 12).Value = 24380  # L66 was 26725
$ws_ARM.Cells.Item(66, 13).Value = -6039.1535  # M66 was -6865.916499999999
$ws_ARM.Cells.Item(66, 14).Value = -31244  # N66 was -33589
$ws_ARM.Cells.Item(74, 8).Value = 1855.48  # H74 was 2067.8572
$ws_ARM.Cells.Item(74, 9).Value = 1575.8334  # I74 was 1778.6
$ws_ARM.Cells.Item(74, 10).Value = 3323.625  # J74 was 3514.1428
$ws_ARM.Cells.Item(74, 11).Value = 1575.8334  # K74 was 1778.6
$ws_ARM.Cells.Item(74, 12).Value = 3323.625  # L74 was 3514.1428
$ws_ARM.Cells.Item(74, 13).Value = -701.8334  # M74 was -904.5999999999999
$ws_ARM.Cells.Item(74, 14).Value = -5071.625  # N74 was -5262.1428
$ws_ARM.Cells.Item(77, 8).Value = 1855.48  # H77 was 2067.8572
$ws_ARM.Cells.Item(77, 9).Value = 1575.8334  # I77 was 1778.6
$ws_ARM.Cells.Item(77, 10).Value = 3323.625  # J77 was 3514.1428
$ws_ARM.Cells.Item(77, 11).Value = 7879.166999999999  # K77 was 8893
$ws_ARM.Cells.Item(77, 12).Value = 16618.125  # L77 was 17570.714
$ws_ARM.Cells.Item(77, 13).Value = -3511.166999999999  # M77 was -4525
$ws_ARM.Cells.Item(77, 14).Value = -25354.125  # N77 was -26306.714
$ws_ARM.Cells.Item(88, 8).Value = 3232146.5  # H88 was 2155288.2
$ws_ARM.Cells.Item(88, 9).Value = 8000  # I88 was 4098.3335
$ws_ARM.Cells.Item(88, 10).Value = 4038183  # J88 was 3589415
$ws_ARM.Cells.Item(88, 11).Value = 8000  # K88 was 4098.3335
$ws_ARM.Cells.Item(88, 12).Value = 4038183  # L88 was 3589415
$ws_ARM.Cells.Item(88, 13).Value = -7594  # M88 was -3692.3335
$ws_ARM.Cells.Item(88, 14).Value = -4038995  # N88 was -3590227
$ws_ARM.Cells.Item(91, 8).Value = 3232146.5  # H91 was 2155288.2
$ws_ARM.Cells.Item(91, 9).Value = 8000  # I91 was 4098.3335
$ws_ARM.Cells.Item(91, 10).Value = 4038183  # J91 was 3589415
$ws_ARM.Cells.Item(91, 11).Value = 8000  # K91 was 4098.3335
$ws_ARM.Cells.Item(91, 12).Value = 4038183  # L91 was 3589415
$ws_ARM.Cells.Item(91, 13).Value = -6596  # M91 was -2694.3335
$ws_ARM.Cells.Item(91, 14).Value = -4040991  # N91 was -3592223
$ws_ARM.Cells.Item(132, 8).Value = 12822861  # H132 was 20836128
$ws_ARM.Cells.Item(132, 9).Value = 15153381  # I132 was 27779836
$ws_ARM.Cells.Item(132, 11).Value = 45460143  # K132 was 83339508
$ws_ARM.Cells.Item(132, 13).Value = -45457613  # M132 was -83336978
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Cells.Item(22, 8).Value = 279.5  # H22 was 0
$ws_BSM.Cells.Item(22, 9).Value = 279.5  # I22 was 0
$ws_BSM.Cells.Item(22, 11).Value = 279.5  # K22 was 0
$ws_BSM.Cells.Item(22, 13).Value = -106.5  # M22 was None
$ws_BSM.Cells.Item(35, 8).Value = 35326.668  # H35 was 35423.332
$ws_BSM.Cells.Item(35, 10).Value = 35326.668  # J35 was 35423.332
$ws_BSM.Cells.Item(35, 12).Value = 35326.668  # L35 was 35423.332
$ws_BSM.Cells.Item(35, 14).Value = -35946.668  # N35 was -36043.332
$ws_BSM.Cells.Item(82, 8).Value = 8628.75  # H82 was 9585.833000000001
$ws_BSM.Cells.Item(82, 9).Value = 4838.3335  # I82 was 4378.75
$ws_BSM.Cells.Item(82, 11).Value = 4838.3335  # K82 was 4378.75
$ws_BSM.Cells.Item(82, 13).Value = -4455.3335  # M82 was -3995.75
$ws_BSM.Cells.Item(85, 8).Value = 8628.75  # H85 was 9585.833000000001
$ws_BSM.Cells.Item(85, 9).Value = 4838.3335  # I85 was 4378.75
$ws_BSM.Cells.Item(85, 11).Value = 4838.3335  # K85 was 4378.75
$ws_BSM.Cells.Item(85, 13).Value = -3512.3335  # M85 was -3052.75
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Cells.Item(19, 8).Value = 868.5714  # H19 was 1073.3334
$ws_CRP.Cells.Item(19, 9).Value = 396.92307  # I19 was 332.5
$ws_CRP.Cells.Item(19, 11).Value = 396.92307  # K19 was 332.5
$ws_CRP.Cells.Item(19, 13).Value = -226.92307  # M19 was -162.5
$ws_CRP.Cells.Item(24, 8).Value = 868.5714  # H24 was 1073.3334
$ws_CRP.Cells.Item(24, 9).Value = 396.92307  # I24 was 332.5
$ws_CRP.Cells.Item(24, 11).Value = 396.92307  # K24 was 332.5
$ws_CRP.Cells.Item(24, 13).Value = -226.92307  # M24 was -162.5
$ws_CRP.Cells.Item(58, 8).Value = 2907.0952  # H58 was 3303.6875
$ws_CRP.Cells.Item(58, 9).Value = 1920.8334  # I58 was 1965
$ws_CRP.Cells.Item(58, 10).Value = 4222.1113  # J58 was 5534.8335
$ws_CRP.Cells.Item(58, 11).Value = 1920.8334  # K58 was 1965
$ws_CRP.Cells.Item(58, 12).Value = 4222.1113  # L58 was 5534.8335
$ws_CRP.Cells.Item(58, 13).Value = -1717.8334  # M58 was -1762
$ws_CRP.Cells.Item(58, 14).Value = -4628.1113  # N58 was -5940.8335
$ws_CRP.Cells.Item(132, 8).Value = 81836.28  # H132 was 64218.566
$ws_CRP.Cells.Item(132, 9).Value = 2149.4443  # I132 was 1756.2727
$ws_CRP.Cells.Item(132, 10).Value = 161523.11  # J132 was 121475.664
$ws_CRP.Cells.Item(132, 11).Value = 6448.3329  # K132 was 5268.8181
$ws_CRP.Cells.Item(132, 12).Value = 484569.33  # L132 was 364426.992
$ws_CRP.Cells.Item(132, 13).Value = -3918.3329  # M132 was -2738.8181
$ws_CRP.Cells.Item(132, 14).Value = -489629.33  # N132 was -369486.992
$ws_CRP.Cells.Item(136, 8).Value = 2907.0952  # H136 was 3303.6875
$ws_CRP.Cells.Item(136, 9).Value = 1920.8334  # I136 was 1965
$ws_CRP.Cells.Item(136, 10).Value = 4222.1113  # J136 was 5534.8335
$ws_CRP.Cells.Item(136, 11).Value = 5762.5002  # K136 was 5895
$ws_CRP.Cells.Item(136, 12).Value = 12666.3339  # L136 was 16604.5005
$ws_CRP.Cells.Item(136, 13).Value = -3212.5002  # M136 was -3345
$ws_CRP.Cells.Item(136, 14).Value = -17766.3339  # N136 was -21704.5005
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Cells.Item(108, 8).Value = 3147.6365  # H108 was 3085.4211
$ws_CUL.Cells.Item(108, 10).Value = 2995.7144  # J108 was 2997.9333
$ws_CUL.Cells.Item(108, 12).Value = 8987.143199999999  # L108 was 8993.7999
$ws_CUL.Cells.Item(108, 14).Value = -14747.1432  # N108 was -14753.7999
$ws_CUL.Cells.Item(113, 8).Value = 2270.8645  # H113 was 2369.125
$ws_CUL.Cells.Item(113, 9).Value = 2929.122  # I113 was 2931.2927
$ws_CUL.Cells.Item(113, 10).Value = 771.5  # J113 was 832.5333000000001
$ws_CUL.Cells.Item(113, 11).Value = 8787.366  # K113 was 8793.8781
$ws_CUL.Cells.Item(113, 12).Value = 2314.5  # L113 was 2497.5999
$ws_CUL.Cells.Item(113, 13).Value = -6617.366  # M113 was -6623.8781
$ws_CUL.Cells.Item(113, 14).Value = -6654.5  # N113 was -6837.5999
$ws_CUL.Cells.Item(116, 8).Value = 1733.3334  # H116 was 1000
$ws_CUL.Cells.Item(116, 9).Value = 2000  # I116 was 1000
$ws_CUL.Cells.Item(116, 10).Value = 1200  # J116 was 0
$ws_CUL.Cells.Item(116, 11).Value = 6000  # K116 was 3000
$ws_CUL.Cells.Item(116, 12).Value = 3600  # L116 was 0
$ws_CUL.Cells.Item(116, 13).Value = -2558  # M116 was 442
$ws_CUL.Cells.Item(116, 14).Value = -10484  # N116 was None
$ws_CUL.Cells.Item(119, 8).Value = 2911.2  # H119 was 2833
$ws_CUL.Cells.Item(119, 9).Value = 1729.6666  # I119 was 1638.2727
$ws_CUL.Cells.Item(119, 10).Value = 4683.5  # J119 was 4710.4287
$ws_CUL.Cells.Item(119, 11).Value = 5188.9998  # K119 was 4914.8181
$ws_CUL.Cells.Item(119, 12).Value = 14050.5  # L119 was 14131.2861
$ws_CUL.Cells.Item(119, 13).Value = -350.9997999999996  # M119 was -76.81810000000041
$ws_CUL.Cells.Item(119, 14).Value = -23726.5  # N119 was -23807.2861
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Cells.Item(102, 8).Value = 1605  # H102 was 1548
$ws_GSM.Cells.Item(102, 9).Value = 1605  # I102 was 1404
$ws_GSM.Cells.Item(102, 10).Value = 0  # J102 was 1980
$ws_GSM.Cells.Item(102, 11).Value = 1605  # K102 was 1404
$ws_GSM.Cells.Item(102, 12).Value = 0  # L102 was 1980
$ws_GSM.Cells.Item(102, 13).Value = 17  # M102 was 218
$ws_GSM.Cells.Item(102, 14).ClearContents()  # N102 was -5224
$ws_GSM.Cells.Item(113, 8).Value = 1622.8182  # H113 was 1954.4286
$ws_GSM.Cells.Item(113, 9).Value = 1685.1  # I113 was 2030.1666
$ws_GSM.Cells.Item(113, 10).Value = 1000  # J113 was 1500
$ws_GSM.Cells.Item(113, 11).Value = 1685.1  # K113 was 2030.1666
$ws_GSM.Cells.Item(113, 12).Value = 1000  # L113 was 1500
$ws_GSM.Cells.Item(113, 13).Value = 484.9000000000001  # M113 was 139.8334
$ws_GSM.Cells.Item(113, 14).Value = -5340  # N113 was -5840
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Cells.Item(2, 8).Value = 70000.8  # H2 was 33319.453
$ws_LTW.Cells.Item(2, 9).Value = 0  # I2 was 1500
$ws_LTW.Cells.Item(2, 10).Value = 70000.8  # J2 was 36501.4
$ws_LTW.Cells.Item(2, 11).Value = 0  # K2 was 1500
$ws_LTW.Cells.Item(2, 12).Value = 70000.8  # L2 was 36501.4
$ws_LTW.Cells.Item(2, 13).ClearContents()  # M2 was -1388
$ws_LTW.Cells.Item(2, 14).Value = -70224.8  # N2 was -36725.4
$ws_LTW.Cells.Item(16, 8).Value = 1757.6285  # H16 was 0
$ws_LTW.Cells.Item(16, 9).Value = 1548.1875  # I16 was 0
$ws_LTW.Cells.Item(16, 10).Value = 3991.6667  # J16 was 0
$ws_LTW.Cells.Item(16, 11).Value = 1548.1875  # K16 was 0
$ws_LTW.Cells.Item(16, 12).Value = 3991.6667  # L16 was 0
$ws_LTW.Cells.Item(16, 13).Value = -1378.1875  # M16 was None
$ws_LTW.Cells.Item(16, 14).Value = -4331.6667  # N16 was None
$ws_LTW.Cells.Item(122, 8).Value = 64864.562  # H122 was 85836.086
$ws_LTW.Cells.Item(122, 9).Value = 93221.17999999999  # I122 was 113514.78
$ws_LTW.Cells.Item(122, 10).Value = 2480  # J122 was 2800
$ws_LTW.Cells.Item(122, 11).Value = 279663.54  # K122 was 340544.34
$ws_LTW.Cells.Item(122, 12).Value = 7440  # L122 was 8400
$ws_LTW.Cells.Item(122, 13).Value = -277213.54  # M122 was -338094.34
$ws_LTW.Cells.Item(122, 14).Value = -12340  # N122 was -13300
$ws_LTW.Cells.Item(132, 8).Value = 2884.8928  # H132 was 2321.9473
$ws_LTW.Cells.Item(132, 9).Value = 2526.7368  # I132 was 2077.3794
$ws_LTW.Cells.Item(132, 10).Value = 3641  # J132 was 3110
$ws_LTW.Cells.Item(132, 11).Value = 7580.2104  # K132 was 6232.138199999999
$ws_LTW.Cells.Item(132, 12).Value = 10923  # L132 was 9330
$ws_LTW.Cells.Item(132, 13).Value = -5050.2104  # M132 was -3702.138199999999
$ws_LTW.Cells.Item(132, 14).Value = -15983  # N132 was -14390
$ws_LTW.Cells.Item(136, 8).Value = 2109.3462  # H136 was 0
$ws_LTW.Cells.Item(136, 9).Value = 1754.619  # I136 was 0
$ws_LTW.Cells.Item(136, 10).Value = 3599.2  # J136 was 0
$ws_LTW.Cells.Item(136, 11).Value = 5263.857  # K136 was 0
$ws_LTW.Cells.Item(136, 12).Value = 10797.6  # L136 was 0
$ws_LTW.Cells.Item(136, 13).Value = -2713.857  # M136 was None
$ws_LTW.Cells.Item(136, 14).Value = -15897.6  # N136 was None
$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Cells.Item(39, 8).Value = 14666.333  # H39 was 9999
$ws_WVR.Cells.Item(39, 9).Value = 9499.5  # I39 was 9999
$ws_WVR.Cells.Item(39, 10).Value = 25000  # J39 was 0
$ws_WVR.Cells.Item(39, 11).Value = 9499.5  # K39 was 9999
$ws_WVR.Cells.Item(39, 12).Value = 25000  # L39 was 0
$ws_WVR.Cells.Item(39, 13).Value = -9086.5  # M39 was -9586
$ws_WVR.Cells.Item(39, 14).Value = -25826  # N39 was None
$ws_WVR.Cells.Item(94, 8).Value = 10000  # H94 was 24250
$ws_WVR.Cells.Item(94, 10).Value = 10000  # J94 was 24250
$ws_WVR.Cells.Item(94, 12).Value = 10000  # L94 was 24250
$ws_WVR.Cells.Item(94, 14).Value = -11802  # N94 was -26052
$ws_WVR.Cells.Item(126, 8).Value = 1731574.1  # H126 was 2102544.5
$ws_WVR.Cells.Item(126, 9).Value = 2451878.2  # I126 was 3269044.2
$ws_WVR.Cells.Item(126, 11).Value = 7355634.600000001  # K126 was 9807132.600000001
$ws_WVR.Cells.Item(126, 13).Value = -7353164.600000001  # M126 was -9804662.600000001
$ws_WVR.Cells.Item(136, 8).Value = 667848.5600000001  # H136 was 687437.3
$ws_WVR.Cells.Item(136, 9).Value = 898328.1  # I136 was 898334.7
$ws_WVR.Cells.Item(136, 10).Value = 2018.5555  # J136 was 2020.875
$ws_WVR.Cells.Item(136, 11).Value = 2694984.3  # K136 was 2695004.1
$ws_WVR.Cells.Item(136, 12).Value = 6055.666499999999  # L136 was 6062.625
$ws_WVR.Cells.Item(136, 13).Value = -2692434.3  # M136 was -2692454.1
$ws_WVR.Cells.Item(136, 14).Value = -11155.6665  # N136 was -11162.625